# Reorders the header columns of the Products, Sales and Expenses sheets
# (Expenses also loses its old "receipt_number" column) and appends a first
# data row to each of those three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Products"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Products")

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "cost_price"
$ws.Range("E1").Value = "category"
$ws.Range("F1").Value = "stock"
$ws.Range("G1").Value = "min_stock"
$ws.Range("H1").Value = "supplier"
$ws.Range("I1").Value = "sku"
$ws.Range("J1").Value = "id"
$ws.Range("K1").Value = "created_date"
$ws.Range("L1").Value = "last_updated"

$ws.Range("A2").Value = "Amul Butter (500g)"
$ws.Range("B2").Value = ""

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1000"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "10"

$ws.Range("E2").Value = "Clothing"
$ws.Range("F2").Value = 450
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "79bf20ff-9535-45ea-9dcd-a9596ae5258a"
$ws.Range("K2").Value = "2025-09-23T11:20:28.437Z"
$ws.Range("L2").Value = "2025-09-23T11:20:38.822Z"

# ---------------------------------------------------------------------
# Sheet "Sales"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sales")

$ws.Range("A1").Value = "product_id"
$ws.Range("B1").Value = "quantity"
$ws.Range("C1").Value = "unit_price"
$ws.Range("D1").Value = "customer_name"
$ws.Range("E1").Value = "payment_method"
$ws.Range("F1").Value = "cashier"
$ws.Range("G1").Value = "notes"
$ws.Range("H1").Value = "id"
$ws.Range("I1").Value = "product_name"
$ws.Range("J1").Value = "total_amount"
$ws.Range("K1").Value = "profit"
$ws.Range("L1").Value = "sale_date"

$ws.Range("A2").Value = "79bf20ff-9535-45ea-9dcd-a9596ae5258a"
$ws.Range("B2").Value = 50

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1000"

$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Cash"
$ws.Range("F2").Value = "Admin"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "cae2620a-9922-4239-8ac5-97c254e4019b"
$ws.Range("I2").Value = "Amul Butter (500g)"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "50000"

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "49500"

$ws.Range("L2").Value = "2025-09-23T11:20:38.815Z"

# ---------------------------------------------------------------------
# Sheet "Expenses" (drops the old "receipt_number" column entirely)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Expenses")

$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "amount"
$ws.Range("D1").Value = "payment_method"
$ws.Range("E1").Value = "vendor"
$ws.Range("F1").Value = "notes"
$ws.Range("G1").Value = "id"
$ws.Range("H1").Value = "expense_date"
$ws.Range("I1").ClearContents()

$ws.Range("A2").Value = "Rent"
$ws.Range("B2").Value = ""

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "5000"

$ws.Range("D2").Value = "Cash"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "ed6c6149-eea7-457f-95a8-eb25cdf5b985"
$ws.Range("H2").Value = "2025-09-23T11:26:47.809Z"
